# Adaptations Data Model for Legal info on Resources
# Adds a new "Authorship Resource" column (M) to Sheet1, populated for every
# data row with the authorship/legal-info string, matching the author cell
# style used elsewhere in the sheet (explicit font color applied, no
# fill/border, default number format).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$header = "Authorship Resource"
$value  = "Noémi Villars-Amberg, Daniela Subotic"

# Header cell
$ws.Range("M1").Value = $header

# Data rows 2..93
for ($r = 2; $r -le 93; $r++) {
    $ws.Cells.Item($r, 13).Value = $value
}

# Apply the distinguishing format (explicit black font color, no fill/border)
# to the whole new column range in one shot so every cell shares a single
# new style entry.
$ws.Range("M1:M93").Font.Color = 0

# Match the post-edit selection recorded in the workbook: the newly added
# column is selected (M2:M93), with M2 as the active cell.
$ws.Range("M2:M93").Select()
